# Adding 16-apr (16 April 2020) age/sex breakdown data to the "Data" sheet.
# Source: tabula-Actualizacion_78_COVID-19(2).csv  (20 rows: 10 age groups x 2 sexes)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$startRow = 511
$dateSerial = 43937          # 2020-04-16
$ageOrder   = @("0-9","10-19","20-29","30-39","40-49","50-59","60-69","70-79","80-89","90 y +")

$rows = @(
    @("Fem","0-9",199,66,0.3,1,0,1,0,0.5),
    @("Fem","10-19",395,80,0.3,3,0.1,1,0,0.3),
    @("Fem","20-29",4501,490,1.9,19,0.8,5,0.1,0.1),
    @("Fem","30-39",7764,1195,4.5999999999999996,49,2,15,0.3,0.2),
    @("Fem","40-49",11221,2341,9,117,4.8,42,0.9,0.4),
    @("Fem","50-59",13684,3810,14.7,243,9.9,103,2.2999999999999998,0.8),
    @("Fem","60-69",10037,4767,18.3,421,17.2,285,6.4,2.8),
    @("Fem","70-79",8847,5781,22.3,695,28.3,864,19.5,9.8000000000000007),
    @("Fem","80-89",10685,5507,21.2,650,26.5,1961,44.3,18.399999999999999),
    @("Fem","90 y +",5467,1944,7.5,256,10.4,1149,26,21),
    @("Masc","0-9",236,93,0.3,13,0.4,0,0,0),
    @("Masc","10-19",345,70,0.2,5,0.1,1,0,0.3),
    @("Masc","20-29",2360,479,1.4,32,0.9,15,0.2,0.6),
    @("Masc","30-39",4903,1337,3.8,114,3.1,22,0.3,0.4),
    @("Masc","40-49",8660,3481,9.8000000000000007,324,8.9,76,1.1000000000000001,0.9),
    @("Masc","50-59",11198,5728,16.2,726,19.899999999999999,216,3.3,1.9),
    @("Masc","60-69",11562,7447,21,1218,33.4,687,10.3,5.9),
    @("Masc","70-79",11810,8845,25,1109,30.4,2062,31.1,17.5),
    @("Masc","80-89",9241,6447,18.2,95,2.6,2698,40.6,29.2),
    @("Masc","90 y +",2667,1461,4.0999999999999996,9,0.2,862,13,32.299999999999997)
)

$r = $startRow
foreach ($row in $rows) {
    $gender = $row[0]
    $ageGroup = $row[1]

    $ws.Range("A$r").Value = $dateSerial
    $ws.Range("A$r").NumberFormat = "d-mmm"

    $ws.Range("B$r").Value = $gender

    $ws.Range("C$r").Value = $ageGroup
    if ($ageGroup -eq "10-19") {
        # Matches the existing column C formatting rule for this label
        $ws.Range("C$r").NumberFormat = "@"
    }

    $cols = @("D","E","F","G","H","I","J","K")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $val = $row[2 + $i]
        $cell = $ws.Range("$col$r")
        $cell.Value = $val
        if ($val -ge 1000) {
            $cell.NumberFormat = "#,##0"
        }
    }

    $r = $r + 1
}

$lastRow = $r - 1

# Register the new tabula range as a worksheet-scoped defined name, same
# naming convention used for the other imported "Actualizacion" ranges.
[void]$ws.Names.Add("tabula_Actualizacion_78_COVID_19_2", $ws.Range("C$startRow`:K$lastRow"))

# Move the selection to the first cell of the newly pasted range, like Excel
# does right after an import/paste operation.
[void]$ws.Range("C$startRow").Select()
